# issue #5: stock data output to json file
#
# The "股票" (stock) sheet gains a new "property_category" column (with
# value "stock" for every data row), inserted between the existing
# "total" and "date" columns. Inserting the column shifts the former
# date / legislator_name / legislator_id columns one slot to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Column H currently holds "date" (G = total, H = date, I = legislator_name,
# J = legislator_id). Insert a new blank column at H, shifting date /
# legislator_name / legislator_id right by one (-> I / J / K).
$ws.Columns.Item(8).Insert()

# Header for the newly inserted column.
$ws.Cells.Item(1, 8).Value = "property_category"

# Every data row on this sheet describes a stock holding.
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 8).Value = "stock"
}
